# Weekly refresh of Jengibre (ginger) price rows: dates and price/volume
# figures get reshuffled across rows 2-27 (rows 5 and 20 are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44453
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14600
$ws.Range("P2").Value = 1123
$ws.Range("D3").Value = 44308
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 26000
$ws.Range("L3").Value = 27000
$ws.Range("M3").Value = 26400
$ws.Range("P3").Value = 2031
$ws.Range("D4").Value = 44474
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("P4").Value = 1038
$ws.Range("D6").Value = 44510
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("P6").Value = 1192
$ws.Range("D7").Value = 44264
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 30000
$ws.Range("L7").Value = 32000
$ws.Range("M7").Value = 31000
$ws.Range("P7").Value = 2385
$ws.Range("D8").Value = 44503
$ws.Range("J8").Value = 35
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15429
$ws.Range("P8").Value = 1187
$ws.Range("D9").Value = 44509
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("P9").Value = 1192
$ws.Range("D10").Value = 44313
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 25000
$ws.Range("L10").Value = 26000
$ws.Range("M10").Value = 25600
$ws.Range("P10").Value = 1969
$ws.Range("D11").Value = 44327
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 24400
$ws.Range("P11").Value = 1877
$ws.Range("D12").Value = 44425
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14500
$ws.Range("P12").Value = 1115
$ws.Range("D13").Value = 44433
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 13500
$ws.Range("P13").Value = 1038
$ws.Range("D14").Value = 44377
$ws.Range("J14").Value = 40
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("P14").Value = 1115
$ws.Range("D15").Value = 44610
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17400
$ws.Range("P15").Value = 1338
$ws.Range("D16").Value = 44350
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 23000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 24000
$ws.Range("P16").Value = 1846
$ws.Range("D17").Value = 44316
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 27000
$ws.Range("L17").Value = 28000
$ws.Range("M17").Value = 27400
$ws.Range("P17").Value = 2108
$ws.Range("D18").Value = 44488
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 16000
$ws.Range("L18").Value = 17000
$ws.Range("M18").Value = 16500
$ws.Range("P18").Value = 1269
$ws.Range("D19").Value = 44334
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 26000
$ws.Range("L19").Value = 28000
$ws.Range("M19").Value = 27200
$ws.Range("P19").Value = 2092
$ws.Range("D21").Value = 44462
$ws.Range("J21").Value = 60
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14500
$ws.Range("P21").Value = 1115
$ws.Range("D22").Value = 44159
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 30000
$ws.Range("L22").Value = 32000
$ws.Range("M22").Value = 31000
$ws.Range("P22").Value = 2385
$ws.Range("D23").Value = 44362
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 15000
$ws.Range("L23").Value = 16000
$ws.Range("M23").Value = 15500
$ws.Range("P23").Value = 1192
$ws.Range("D24").Value = 44355
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 18000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 19000
$ws.Range("P24").Value = 1462
$ws.Range("D25").Value = 44435
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 13000
$ws.Range("L25").Value = 14000
$ws.Range("M25").Value = 13500
$ws.Range("P25").Value = 1038
$ws.Range("D26").Value = 44467
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 13000
$ws.Range("L26").Value = 14000
$ws.Range("M26").Value = 13500
$ws.Range("P26").Value = 1038
$ws.Range("D27").Value = 44383
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 15000
$ws.Range("L27").Value = 16000
$ws.Range("M27").Value = 15400
$ws.Range("P27").Value = 1185
